$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.413.40'
$ws.Range('E2').Value = '  +0.08%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.850.18'
$ws.Range('E3').Value = '  +0.21%  '

$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('E5').Value = '  +0.26%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6274'
$ws.Range('E6').Value = '  -0.39%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.002'
$ws.Range('E7').Value = '  +0.14%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07634'
$ws.Range('E8').Value = '  +0.56%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2908'
$ws.Range('E9').Value = '  -0.68%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.74'
$ws.Range('E10').Value = '  +1.03%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07748'
$ws.Range('E11').Value = '  +0.18%  '

$ws.Range('E12').Value = '  +0.70%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6788'
$ws.Range('E13').Value = '  +0.07%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.00001064'
$ws.Range('E14').Value = '  -2.38%  '

$ws.Range('E15').Value = '  -0.49%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.160'
$ws.Range('E16').Value = '  +0.13%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.423.03'
$ws.Range('E17').Value = '  +0.04%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '226.45'
$ws.Range('E18').Value = '  -0.90%  '

$ws.Range('E19').Value = '  -0.85%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.10%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.476'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.14%  '

$ws.Range('E23').Value = '  +0.32%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.420'
$ws.Range('E25').Value = '  +0.63%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.69'
$ws.Range('E26').Value = '  +0.46%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.384'
$ws.Range('E27').Value = '  +6.43%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.467'
$ws.Range('E28').Value = '  +0.20%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05588'
$ws.Range('E29').Value = '  -0.26%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.128'
$ws.Range('E30').Value = '  +0.63%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.055'
$ws.Range('E31').Value = '  +0.51%  '

$ws.Range('E32').Value = '  -0.38%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.6958'
$ws.Range('E34').Value = '  -1.81%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.593'
$ws.Range('E35').Value = '  +0.29%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.233.41'
$ws.Range('E36').Value = '  +0.30%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01802'
$ws.Range('E37').Value = '  +0.30%  '

$ws.Range('E38').Value = '  -1.68%  '

$ws.Range('E39').Value = '  -0.47%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9035'
$ws.Range('E40').Value = '  -0.45%  '

$ws.Range('E41').Value = '  +0.20%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '101.62'
$ws.Range('E42').Value = '  -0.22%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '65.92'
$ws.Range('E43').Value = '  -0.14%  '

$ws.Range('E44').Value = '  -0.83%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.172'
$ws.Range('E45').Value = '  -0.07%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4012'
$ws.Range('E46').Value = '  -0.12%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.978'
$ws.Range('E47').Value = '  +0.20%  '

$ws.Range('E48').Value = '  -0.07%  '

$ws.Range('E49').Value = '  +2.02%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05706'
$ws.Range('E50').Value = '  -0.02%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4635'
$ws.Range('E51').Value = '  +0.24%  '
